# Weekly refresh of the Jengibre (ginger) price sheet.
#
# The underlying rows (3-11, 13-20; row 12 is untouched) simply get
# reshuffled / re-dated, i.e. each row's Fecha/Calidad/Volumen/Precio
# minimo/maximo/promedio/Precio-$-Kg set gets replaced with the values
# that belong to the updated weekly snapshot. All other columns
# (Mercado ID, Mercado, Region, Codreg, Categoria, Variedad, Unidad de
# comercializacion, Origen, Kg o Unidades, Clasificacion) stay the same
# across every row, so only D/I/J/K/L/M/P need to be rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44918
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 12000
$ws.Range("L3").Value = 13000
$ws.Range("M3").Value = 12750
$ws.Range("P3").Value = 981

$ws.Range("D4").Value = 44943
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 14429
$ws.Range("P4").Value = 1110

$ws.Range("D5").Value = 44159
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 23000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 23500
$ws.Range("P5").Value = 1808

$ws.Range("D6").Value = 44893
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 900
$ws.Range("K6").Value = 13000
$ws.Range("L6").Value = 14000
$ws.Range("M6").Value = 13444
$ws.Range("P6").Value = 1034

$ws.Range("D7").Value = 44406
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("P7").Value = 1346

$ws.Range("D8").Value = 44890
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("P8").Value = 1115

$ws.Range("D9").Value = 44389
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 12500
$ws.Range("P9").Value = 962

$ws.Range("D10").Value = 44469
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13500
$ws.Range("P10").Value = 1038

$ws.Range("D11").Value = 44580
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 11500
$ws.Range("P11").Value = 885

$ws.Range("D13").Value = 44397
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 12500
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12750
$ws.Range("P13").Value = 981

$ws.Range("D14").Value = 44379
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 12667
$ws.Range("P14").Value = 974

$ws.Range("D15").Value = 44914
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14400
$ws.Range("P15").Value = 1108

$ws.Range("D16").Value = 44832
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 13000
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = 13500
$ws.Range("P16").Value = 1038

$ws.Range("D17").Value = 44855
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = 10000
$ws.Range("P17").Value = 769

$ws.Range("D18").Value = 44616
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 19000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 19500
$ws.Range("P18").Value = 1500

$ws.Range("D19").Value = 44320
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 160
$ws.Range("K19").Value = 19000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19500
$ws.Range("P19").Value = 1500

$ws.Range("D20").Value = 44592
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 13000
$ws.Range("M20").Value = 12500
$ws.Range("P20").Value = 962
